$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 45
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03)
$ws.Range("C2:C45").Value = 45202
